$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "F:/My Files" path prefix with "D:/My Files" in column A (rows 2-9)
for ($r = 2; $r -le 9; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    if ($old -and $old.StartsWith("F:/My Files")) {
        $cell.Value = $old.Replace("F:/My Files", "D:/My Files")
    }
}

# Update selection to A2
$ws.Range("A2").Select()
